# Update version/build string across workbook for the
# "Coal Mine Boundaries and Methane Sources - version 1.0.0" release.

$wb = $excel.ActiveWorkbook

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = "Version: " + $newVersion

$newCitation = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Shanxi Jinyuan Coal Mine, China, M0330, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"
$wsAbout.Range("A6").Value = $newCitation

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 14; $row++) {
    $wsData.Cells.Item($row, 19).Value = $newVersion
}
